$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$sheetInfo = @(
    @{ Name = "zh-cn"; HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a80beaea57bc5b894c612f7f7c63c0788b7e633f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"; HandbackTime2 = "2016-02-24 12:01:03" },
    @{ Name = "de-de"; HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a059014b14679a163fb395b5dff9bb4ff9fd7b57/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"; HandbackTime2 = "2016-02-24 12:01:33" }
)

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/df479055007b6b76fd6e4bd3b5f5ff0baac197a5/e2e"

foreach ($info in $sheetInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Row 2: 0203c660-bb9c-47d6-b768-b5c3eddf0d26
    $ws.Range("B2").Value = $newStatus

    $ws.Range("E2").Value = $ws.Range("A2").Value
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("E2"), ($mdBase + "/0203c660-bb9c-47d6-b768-b5c3eddf0d26.md"), [Type]::Missing, [Type]::Missing, "0203c660-bb9c-47d6-b768-b5c3eddf0d26.md") | Out-Null
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = 15570276

    $ws.Range("F2").Value = $ws.Range("C2").Value
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("F2"), ($info.HandoffBase + "/0203c660-bb9c-47d6-b768-b5c3eddf0d26.a878ecf2046da54bb8adbb77b490e7adb024db43." + $info.Name + ".xlf"), [Type]::Missing, [Type]::Missing, ("0203c660-bb9c-47d6-b768-b5c3eddf0d26.a878ecf2046da54bb8adbb77b490e7adb024db43." + $info.Name + ".xlf")) | Out-Null
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 15570276

    $ws.Range("G2").Value = $info.HandbackTime2

    # Row 3: 95e730fc-4864-4976-9d13-5bca16481286
    $ws.Range("B3").Value = $newStatus

    $ws.Range("E3").Value = $ws.Range("A3").Value
    $ws.Range("E3").Font.Underline = $true
    $ws.Range("E3").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("E3"), ($mdBase + "/95e730fc-4864-4976-9d13-5bca16481286.md"), [Type]::Missing, [Type]::Missing, "95e730fc-4864-4976-9d13-5bca16481286.md") | Out-Null
    $ws.Range("E3").Font.Underline = $true
    $ws.Range("E3").Font.Color = 15570276

    $ws.Range("F3").Value = $ws.Range("C3").Value
    $ws.Range("F3").Font.Underline = $true
    $ws.Range("F3").Font.Color = 15570276
    $ws.Hyperlinks.Add($ws.Range("F3"), ($info.HandoffBase + "/95e730fc-4864-4976-9d13-5bca16481286.498eb5a6334eae58c1726356a086d00994128c4d." + $info.Name + ".xlf"), [Type]::Missing, [Type]::Missing, ("95e730fc-4864-4976-9d13-5bca16481286.498eb5a6334eae58c1726356a086d00994128c4d." + $info.Name + ".xlf")) | Out-Null
    $ws.Range("F3").Font.Underline = $true
    $ws.Range("F3").Font.Color = 15570276

    $ws.Range("G3").Value = $info.HandbackTime2
}
